$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1:E7").Style = "Normal"

$ws.Range("E1").Value = " Oct 07"
$ws.Range("E2").Value = 27
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
